$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.611874666666667
$ws.Range("H2").Value = 4.835624
$ws.Range("I2").Value = 0.06646895152072402
$ws.Range("J2").Value = 0.06646895152072402
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 5.571066532078222
$ws.Range("R2").Value = 50.139598788704
$ws.Range("S2").Value = 0.0006541565664899854
$ws.Range("T2").Value = 0.0006541565664899853

# Row 3
$ws.Range("G3").Value = 1.611874666666667
$ws.Range("H3").Value = 4.835624
$ws.Range("I3").Value = 0.06646895152072402
$ws.Range("J3").Value = 0.06646895152072402
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 486.1416390987005
$ws.Range("R3").Value = 4375.274751888304
$ws.Range("S3").Value = 0.05708292005300976
$ws.Range("T3").Value = 0.05708292005300975

# Row 4
$ws.Range("G4").Value = 1.611874666666667
$ws.Range("H4").Value = 4.835624
$ws.Range("I4").Value = 0.06646895152072402
$ws.Range("J4").Value = 0.06646895152072402
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 74.3642401780418
$ws.Range("R4").Value = 669.2781616023761
$ws.Range("S4").Value = 0.008731874901224283
$ws.Range("T4").Value = 0.008731874901224282

# Row 5
$ws.Range("I5").Value = 0.6469909869698216
$ws.Range("J5").Value = 0.6469909869698216
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 54.22727080236888
$ws.Range("R5").Value = 488.04543722132
$ws.Range("S5").Value = 0.006367384965508107
$ws.Range("T5").Value = 0.006367384965508106

# Row 6
$ws.Range("I6").Value = 0.6469909869698216
$ws.Range("J6").Value = 0.6469909869698216
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.555629868370969
$ws.Range("T6").Value = 0.555629868370969

# Row 7
$ws.Range("I7").Value = 0.6469909869698216
$ws.Range("J7").Value = 0.6469909869698216
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.0849937336333446
$ws.Range("T7").Value = 0.08499373363334459

# Row 8
$ws.Range("I8").Value = 0.2865400615094543
$ws.Range("J8").Value = 0.2865400615094543
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 24.01623179323422
$ws.Range("R8").Value = 216.146086139108
$ws.Range("S8").Value = 0.002819994275679409
$ws.Range("T8").Value = 0.002819994275679408

# Row 9
$ws.Range("I9").Value = 0.2865400615094543
$ws.Range("J9").Value = 0.2865400615094543
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.2460779514180987
$ws.Range("T9").Value = 0.2460779514180986

# Row 10
$ws.Range("I10").Value = 0.2865400615094543
$ws.Range("J10").Value = 0.2865400615094543
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.03764211581567629
$ws.Range("T10").Value = 0.03764211581567629
